$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.74211658298923
$ws.Range("C2").Value = 9.106340609433479
$ws.Range("E2").Value = 10.777995426789
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 33.71702522009998
$ws.Range("H2").Value = 15.2291425669893
$ws.Range("I2").Value = 22.95232803070489
$ws.Range("L2").Value = 9.971151084086252
$ws.Range("N2").Value = 17.38348616892499
$ws.Range("B3").Value = 16.22601943643508
$ws.Range("C3").Value = 8.731812757914806
$ws.Range("E3").Value = 10.80581040853205
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 33.57042645172668
$ws.Range("H3").Value = 15.26946433217166
$ws.Range("I3").Value = 23.04495424042743
$ws.Range("L3").Value = 9.947488199921247
$ws.Range("N3").Value = 17.44783346892727
$ws.Range("B4").Value = 15.90400146914756
$ws.Range("C4").Value = 8.491610540611839
$ws.Range("E4").Value = 10.82428308507186
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 33.49488681808001
$ws.Range("H4").Value = 15.29795471072259
$ws.Range("I4").Value = 23.10800238295127
$ws.Range("L4").Value = 9.934883713706315
$ws.Range("N4").Value = 17.48930575903287
$ws.Range("B5").Value = 15.77169840435675
$ws.Range("C5").Value = 8.391228676775132
$ws.Range("E5").Value = 10.83216165746712
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 33.46776099375218
$ws.Range("H5").Value = 15.3104999451569
$ws.Range("I5").Value = 23.13524078425337
$ws.Range("L5").Value = 9.930234732055579
$ws.Range("N5").Value = 17.50670105135336
$ws.Range("B6").Value = 15.74967095967739
$ws.Range("C6").Value = 8.374411988193517
$ws.Range("E6").Value = 10.83349108311013
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 33.46347804470602
$ws.Range("H6").Value = 15.31263945026449
$ws.Range("I6").Value = 23.13985685261571
$ws.Range("L6").Value = 9.929492306345679
$ws.Range("N6").Value = 17.50961946551247
$ws.Range("B7").Value = 15.90222125760373
$ws.Range("C7").Value = 8.490266757299414
$ws.Range("E7").Value = 10.82438791753696
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 33.49450616346
$ws.Range("H7").Value = 15.29812011836837
$ws.Range("I7").Value = 23.10836348025032
$ws.Range("L7").Value = 9.934819037967671
$ws.Range("N7").Value = 17.48953835191023
$ws.Range("B8").Value = 16.56535430921319
$ws.Range("C8").Value = 8.979372006814659
$ws.Range("E8").Value = 10.78729681388563
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 33.66348704777133
$ws.Range("H8").Value = 15.24226878650783
$ws.Range("I8").Value = 22.98297966126505
$ws.Range("L8").Value = 9.962594885040534
$ws.Range("N8").Value = 17.40526660720551
$ws.Range("B9").Value = 17.81647170848499
$ws.Range("C9").Value = 9.854379739930245
$ws.Range("E9").Value = 10.72561302011153
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 34.10857024227971
$ws.Range("H9").Value = 15.16250809550298
$ws.Range("I9").Value = 22.78641139040166
$ws.Range("L9").Value = 10.03216293625049
$ws.Range("N9").Value = 17.25551799346981
$ws.Range("B10").Value = 18.69526437388099
$ws.Range("C10").Value = 10.44278109550958
$ws.Range("E10").Value = 10.68701732351678
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 34.50292691644402
$ws.Range("H10").Value = 15.12224289119255
$ws.Range("I10").Value = 22.67247139740244
$ws.Range("L10").Value = 10.09222448274447
$ws.Range("N10").Value = 17.15485942200536
$ws.Range("B11").Value = 19.08440244778794
$ws.Range("C11").Value = 10.69812406366997
$ws.Range("E11").Value = 10.6709161016268
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 34.69640220896814
$ws.Range("H11").Value = 15.10794307393456
$ws.Range("I11").Value = 22.62734226609771
$ws.Range("L11").Value = 10.12142854203831
$ws.Range("N11").Value = 17.11108058649051
$ws.Range("B12").Value = 19.23009488584361
$ws.Range("C12").Value = 10.79300992262638
$ws.Range("E12").Value = 10.66502816155646
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 34.77163581708352
$ws.Range("H12").Value = 15.10310837089898
$ws.Range("I12").Value = 22.61122380645803
$ws.Range("L12").Value = 10.13275199347254
$ws.Range("N12").Value = 17.09479046792849
$ws.Range("B13").Value = 19.19879352680911
$ws.Range("C13").Value = 10.77265541429688
$ws.Range("E13").Value = 10.66628693105532
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 34.75534630222869
$ws.Range("H13").Value = 15.10412376033754
$ws.Range("I13").Value = 22.61465190998129
$ws.Range("L13").Value = 10.13030162239154
$ws.Range("N13").Value = 17.09828604809067
$ws.Range("B14").Value = 19.09642277021218
$ws.Range("C14").Value = 10.7059667738668
$ws.Range("E14").Value = 10.67042750457074
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 34.70255260685057
$ws.Range("H14").Value = 15.10753367487718
$ws.Range("I14").Value = 22.62599669108928
$ws.Range("L14").Value = 10.12235486424938
$ws.Range("N14").Value = 17.10973462417741
$ws.Range("B15").Value = 19.03349697638326
$ws.Range("C15").Value = 10.66488179950358
$ws.Range("E15").Value = 10.6729909705774
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 34.67046957565883
$ws.Range("H15").Value = 15.10969799724232
$ws.Range("I15").Value = 22.63307235312471
$ws.Range("L15").Value = 10.11752150519074
$ws.Range("N15").Value = 17.11678467422109
$ws.Range("B16").Value = 18.66960733418707
$ws.Range("C16").Value = 10.42584273264341
$ws.Range("E16").Value = 10.68809886498632
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 34.49056151741696
$ws.Range("H16").Value = 15.1232584856769
$ws.Range("I16").Value = 22.67555618118273
$ws.Range("L16").Value = 10.09035328285087
$ws.Range("N16").Value = 17.15776083694682
$ws.Range("B17").Value = 18.44354813247856
$ws.Range("C17").Value = 10.27601689711443
$ws.Range("E17").Value = 10.69773991852594
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 34.38376193644922
$ws.Range("H17").Value = 15.13260824912641
$ws.Range("I17").Value = 22.70334048894421
$ws.Range("L17").Value = 10.07416421135339
$ws.Range("N17").Value = 17.18341262668498
$ws.Range("B18").Value = 18.31253164549742
$ws.Range("C18").Value = 10.18868274572366
$ws.Range("E18").Value = 10.70342226669274
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 34.32366267297201
$ws.Range("H18").Value = 15.13836394165452
$ws.Range("I18").Value = 22.71995176721005
$ws.Range("L18").Value = 10.06503019453746
$ws.Range("N18").Value = 17.19835623893544
$ws.Range("B19").Value = 18.2680055511338
$ws.Range("C19").Value = 10.15891512411821
$ws.Range("E19").Value = 10.70536975845374
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 34.30354397694712
$ws.Range("H19").Value = 15.14037754885227
$ws.Range("I19").Value = 22.72568412624778
$ws.Range("L19").Value = 10.06196824002337
$ws.Range("N19").Value = 17.20344845386506
$ws.Range("B20").Value = 18.46771635471218
$ws.Range("C20").Value = 10.29208626321754
$ws.Range("E20").Value = 10.6966994268607
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 34.39499377344659
$ws.Range("H20").Value = 15.13157381357468
$ws.Range("I20").Value = 22.70031749750536
$ws.Range("L20").Value = 10.07586923704969
$ws.Range("N20").Value = 17.18066235865076
$ws.Range("B21").Value = 19.12653777166113
$ws.Range("C21").Value = 10.72560413388089
$ws.Range("E21").Value = 10.66920564016169
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 34.7180064371825
$ws.Range("H21").Value = 15.10651633021792
$ws.Range("I21").Value = 22.62263804643729
$ws.Range("L21").Value = 10.12468189142245
$ws.Range("N21").Value = 17.10636409403146
$ws.Range("B22").Value = 19.5473530388091
$ws.Range("C22").Value = 10.99838844980966
$ws.Range("E22").Value = 10.65245642024238
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 34.94055947524072
$ws.Range("H22").Value = 15.09352347941824
$ws.Range("I22").Value = 22.5775327148722
$ws.Range("L22").Value = 10.15812245762695
$ws.Range("N22").Value = 17.05948387058161
$ws.Range("B23").Value = 19.3236909271611
$ws.Range("C23").Value = 10.8537732031994
$ws.Range("E23").Value = 10.66128425428266
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 34.82075130086524
$ws.Range("H23").Value = 15.10014757989925
$ws.Range("I23").Value = 22.60108584203672
$ws.Range("L23").Value = 10.14013591475833
$ws.Range("N23").Value = 17.08435159577239
$ws.Range("B24").Value = 18.45679316769367
$ws.Range("C24").Value = 10.2848250295052
$ws.Range("E24").Value = 10.69716939836725
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 34.38991180238039
$ws.Range("H24").Value = 15.13204029713337
$ws.Range("I24").Value = 22.70168220595007
$ws.Range("L24").Value = 10.07509785467744
$ws.Range("N24").Value = 17.18190514398222
$ws.Range("B25").Value = 17.48441862432543
$ws.Range("C25").Value = 9.627034761519896
$ws.Range("E25").Value = 10.7411184069204
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 33.976167185546
$ws.Range("H25").Value = 15.18087862260726
$ws.Range("I25").Value = 22.83426568857381
$ws.Range("L25").Value = 10.01175109485537
$ws.Range("N25").Value = 17.29437835276153

Write-Output "Applied 216 cell updates to loading_percent sheet (Case_2_111, 380 kV)"
